# Auto-generated edit script applying scheduled-runner price updates
# to the Masamune_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(108, 8).Value = 36625   # H108: 37129.5 -> 36625
$ws.Cells.Item(108, 10).Value = 36625   # J108: 37129.5 -> 36625
$ws.Cells.Item(108, 12).Value = 36625   # L108: 37129.5 -> 36625
$ws.Cells.Item(108, 14).Value = -44305   # N108: -44809.5 -> -44305

$ws.Cells.Item(120, 8).Value = 47711.6   # H120: 47713.2 -> 47711.6
$ws.Cells.Item(120, 10).Value = 47711.6   # J120: 47713.2 -> 47711.6
$ws.Cells.Item(120, 12).Value = 47711.6   # L120: 47713.2 -> 47711.6
$ws.Cells.Item(120, 14).Value = -57387.6   # N120: -57389.2 -> -57387.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 52103.332   # H80: 48100.668 -> 52103.332
$ws.Cells.Item(80, 10).Value = 52103.332   # J80: 48100.668 -> 52103.332
$ws.Cells.Item(80, 12).Value = 52103.332   # L80: 48100.668 -> 52103.332
$ws.Cells.Item(80, 14).Value = -54099.332   # N80: -50096.668 -> -54099.332

$ws.Cells.Item(83, 8).Value = 52103.332   # H83: 48100.668 -> 52103.332
$ws.Cells.Item(83, 10).Value = 52103.332   # J83: 48100.668 -> 52103.332
$ws.Cells.Item(83, 12).Value = 156309.996   # L83: 144302.004 -> 156309.996
$ws.Cells.Item(83, 14).Value = -166293.996   # N83: -154286.004 -> -166293.996

$ws.Cells.Item(107, 8).Value = 42114   # H107: 39996 -> 42114
$ws.Cells.Item(107, 10).Value = 42114   # J107: 39996 -> 42114
$ws.Cells.Item(107, 12).Value = 42114   # L107: 39996 -> 42114
$ws.Cells.Item(107, 14).Value = -49794   # N107: -47676 -> -49794

$ws.Cells.Item(109, 8).Value = 47377   # H109: 45041 -> 47377
$ws.Cells.Item(109, 10).Value = 47377   # J109: 45041 -> 47377
$ws.Cells.Item(109, 12).Value = 47377   # L109: 45041 -> 47377
$ws.Cells.Item(109, 14).Value = -50151   # N109: -47815 -> -50151

$ws.Cells.Item(111, 8).Value = 48986.668   # H111: 48992 -> 48986.668
$ws.Cells.Item(111, 10).Value = 48986.668   # J111: 48992 -> 48986.668
$ws.Cells.Item(111, 12).Value = 48986.668   # L111: 48992 -> 48986.668
$ws.Cells.Item(111, 14).Value = -57166.668   # N111: -57172 -> -57166.668

$ws.Cells.Item(120, 8).Value = 43825.5   # H120: 43575.5 -> 43825.5
$ws.Cells.Item(120, 10).Value = 43825.5   # J120: 43575.5 -> 43825.5
$ws.Cells.Item(120, 12).Value = 43825.5   # L120: 43575.5 -> 43825.5
$ws.Cells.Item(120, 14).Value = -53501.5   # N120: -53251.5 -> -53501.5

$ws.Cells.Item(135, 8).Value = 60465.6   # H135: 48459.668 -> 60465.6
$ws.Cells.Item(135, 10).Value = 60465.6   # J135: 48459.668 -> 60465.6
$ws.Cells.Item(135, 12).Value = 60465.6   # L135: 48459.668 -> 60465.6
$ws.Cells.Item(135, 14).Value = -70605.60000000001   # N135: -58599.668 -> -70605.60000000001

$ws.Cells.Item(139, 8).Value = 40365.8   # H139: 35596.582 -> 40365.8
$ws.Cells.Item(139, 10).Value = 40365.8   # J139: 35596.582 -> 40365.8
$ws.Cells.Item(139, 12).Value = 40365.8   # L139: 35596.582 -> 40365.8
$ws.Cells.Item(139, 14).Value = -50645.8   # N139: -45876.582 -> -50645.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 10475   # H81: 9866.666999999999 -> 10475
$ws.Cells.Item(81, 10).Value = 10475   # J81: 9866.666999999999 -> 10475
$ws.Cells.Item(81, 12).Value = 10475   # L81: 9866.666999999999 -> 10475
$ws.Cells.Item(81, 14).Value = -12597   # N81: -11988.667 -> -12597

$ws.Cells.Item(84, 8).Value = 10475   # H84: 9866.666999999999 -> 10475
$ws.Cells.Item(84, 10).Value = 10475   # J84: 9866.666999999999 -> 10475
$ws.Cells.Item(84, 12).Value = 31425   # L84: 29600.001 -> 31425
$ws.Cells.Item(84, 14).Value = -42033   # N84: -40208.001 -> -42033

$ws.Cells.Item(94, 8).Value = 616.2941   # H94: 615.2 -> 616.2941
$ws.Cells.Item(94, 9).Value = 479.72726   # I94: 482.8 -> 479.72726
$ws.Cells.Item(94, 10).Value = 866.6667   # J94: 880 -> 866.6667
$ws.Cells.Item(94, 11).Value = 479.72726   # K94: 482.8 -> 479.72726
$ws.Cells.Item(94, 12).Value = 866.6667   # L94: 880 -> 866.6667
$ws.Cells.Item(94, 13).Value = -28.72726   # M94: -31.80000000000001 -> -28.72726
$ws.Cells.Item(94, 14).Value = -1768.6667   # N94: -1782 -> -1768.6667

$ws.Cells.Item(110, 8).Value = 48694   # H110: 48702 -> 48694
$ws.Cells.Item(110, 10).Value = 48694   # J110: 48702 -> 48694
$ws.Cells.Item(110, 12).Value = 48694   # L110: 48702 -> 48694
$ws.Cells.Item(110, 14).Value = -56874   # N110: -56882 -> -56874

$ws.Cells.Item(117, 8).Value = 49911.332   # H117: 49914 -> 49911.332
$ws.Cells.Item(117, 10).Value = 49911.332   # J117: 49914 -> 49911.332
$ws.Cells.Item(117, 12).Value = 49911.332   # L117: 49914 -> 49911.332
$ws.Cells.Item(117, 14).Value = -59089.332   # N117: -59092 -> -59089.332

$ws.Cells.Item(119, 8).Value = 47425   # H119: 47091.668 -> 47425
$ws.Cells.Item(119, 10).Value = 47425   # J119: 47091.668 -> 47425
$ws.Cells.Item(119, 12).Value = 47425   # L119: 47091.668 -> 47425
$ws.Cells.Item(119, 14).Value = -57101   # N119: -56767.668 -> -57101

$ws.Cells.Item(120, 8).Value = 48753   # H120: 48761 -> 48753
$ws.Cells.Item(120, 10).Value = 48753   # J120: 48761 -> 48753
$ws.Cells.Item(120, 12).Value = 48753   # L120: 48761 -> 48753
$ws.Cells.Item(120, 14).Value = -58429   # N120: -58437 -> -58429

$ws.Cells.Item(132, 8).Value = 51197.5   # H132: 51826.25 -> 51197.5
$ws.Cells.Item(132, 10).Value = 51197.5   # J132: 51826.25 -> 51197.5
$ws.Cells.Item(132, 12).Value = 51197.5   # L132: 51826.25 -> 51197.5
$ws.Cells.Item(132, 14).Value = -61317.5   # N132: -61946.25 -> -61317.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 500   # H3: 0 -> 500
$ws.Cells.Item(3, 9).Value = 500   # I3: 0 -> 500
$ws.Cells.Item(3, 11).Value = 500   # K3: 0 -> 500
$ws.Cells.Item(3, 13).Value = -387   # M3: None -> -387

$ws.Cells.Item(4, 8).Value = 47814050   # H4: 52846580 -> 47814050
$ws.Cells.Item(4, 10).Value = 204749.1   # J4: 226943.56 -> 204749.1
$ws.Cells.Item(4, 12).Value = 204749.1   # L4: 226943.56 -> 204749.1
$ws.Cells.Item(4, 14).Value = -204973.1   # N4: -227167.56 -> -204973.1

$ws.Cells.Item(116, 8).Value = 47814.332   # H116: 47819.668 -> 47814.332
$ws.Cells.Item(116, 10).Value = 47814.332   # J116: 47819.668 -> 47814.332
$ws.Cells.Item(116, 12).Value = 47814.332   # L116: 47819.668 -> 47814.332
$ws.Cells.Item(116, 14).Value = -56992.332   # N116: -56997.668 -> -56992.332

$ws.Cells.Item(118, 8).Value = 44734   # H118: 44742 -> 44734
$ws.Cells.Item(118, 10).Value = 44734   # J118: 44742 -> 44734
$ws.Cells.Item(118, 12).Value = 44734   # L118: 44742 -> 44734
$ws.Cells.Item(118, 14).Value = -48048   # N118: -48056 -> -48048

$ws.Cells.Item(134, 8).Value = 2348.9412   # H134: 2538.6428 -> 2348.9412
$ws.Cells.Item(134, 9).Value = 1288.2222   # I134: 1389.1111 -> 1288.2222
$ws.Cells.Item(134, 10).Value = 3542.25   # J134: 4607.8 -> 3542.25
$ws.Cells.Item(134, 11).Value = 3864.6666   # K134: 4167.3333 -> 3864.6666
$ws.Cells.Item(134, 12).Value = 10626.75   # L134: 13823.4 -> 10626.75
$ws.Cells.Item(134, 13).Value = -1329.6666   # M134: -1632.3333 -> -1329.6666
$ws.Cells.Item(134, 14).Value = -15696.75   # N134: -18893.4 -> -15696.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 44888.75   # H131: 33387.152 -> 44888.75
$ws.Cells.Item(131, 9).Value = 14883.857   # I131: 17284.5 -> 14883.857
$ws.Cells.Item(131, 10).Value = 50011.535   # J131: 35024.71 -> 50011.535
$ws.Cells.Item(131, 11).Value = 44651.571   # K131: 51853.5 -> 44651.571
$ws.Cells.Item(131, 12).Value = 150034.605   # L131: 105074.13 -> 150034.605
$ws.Cells.Item(131, 13).Value = -39611.571   # M131: -46813.5 -> -39611.571
$ws.Cells.Item(131, 14).Value = -160114.605   # N131: -115154.13 -> -160114.605

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 1160   # H13: 1501.25 -> 1160
$ws.Cells.Item(13, 9).Value = 1160   # I13: 1501.25 -> 1160
$ws.Cells.Item(13, 11).Value = 1160   # K13: 1501.25 -> 1160
$ws.Cells.Item(13, 13).Value = -1021   # M13: -1362.25 -> -1021

$ws.Cells.Item(110, 8).Value = 30996   # H110: 37564.668 -> 30996
$ws.Cells.Item(110, 10).Value = 30996   # J110: 37564.668 -> 30996
$ws.Cells.Item(110, 12).Value = 30996   # L110: 37564.668 -> 30996
$ws.Cells.Item(110, 14).Value = -39176   # N110: -45744.668 -> -39176

$ws.Cells.Item(116, 8).Value = 49734   # H116: 49742 -> 49734
$ws.Cells.Item(116, 10).Value = 49734   # J116: 49742 -> 49734
$ws.Cells.Item(116, 12).Value = 49734   # L116: 49742 -> 49734
$ws.Cells.Item(116, 14).Value = -58912   # N116: -58920 -> -58912

$ws.Cells.Item(122, 8).Value = 3466.6667   # H122: 2850 -> 3466.6667
$ws.Cells.Item(122, 9).Value = 3533.3333   # I122: 2100 -> 3533.3333
$ws.Cells.Item(122, 10).Value = 3400   # J122: 3600 -> 3400
$ws.Cells.Item(122, 11).Value = 10599.9999   # K122: 6300 -> 10599.9999
$ws.Cells.Item(122, 12).Value = 10200   # L122: 10800 -> 10200
$ws.Cells.Item(122, 13).Value = -8149.999899999999   # M122: -3850 -> -8149.999899999999
$ws.Cells.Item(122, 14).Value = -15100   # N122: -15700 -> -15100

$ws.Cells.Item(141, 8).Value = 55685.4   # H141: 62309.5 -> 55685.4
$ws.Cells.Item(141, 10).Value = 55685.4   # J141: 62309.5 -> 55685.4
$ws.Cells.Item(141, 12).Value = 55685.4   # L141: 62309.5 -> 55685.4
$ws.Cells.Item(141, 14).Value = -66045.39999999999   # N141: -72669.5 -> -66045.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(111, 8).Value = 36122.668   # H111: 36443.168 -> 36122.668
$ws.Cells.Item(111, 10).Value = 36122.668   # J111: 36443.168 -> 36122.668
$ws.Cells.Item(111, 12).Value = 36122.668   # L111: 36443.168 -> 36122.668
$ws.Cells.Item(111, 14).Value = -44302.668   # N111: -44623.168 -> -44302.668

$ws.Cells.Item(112, 8).Value = 36902.668   # H112: 40250.668 -> 36902.668
$ws.Cells.Item(112, 10).Value = 36902.668   # J112: 40250.668 -> 36902.668
$ws.Cells.Item(112, 12).Value = 36902.668   # L112: 40250.668 -> 36902.668
$ws.Cells.Item(112, 14).Value = -39856.668   # N112: -43204.668 -> -39856.668

$ws.Cells.Item(132, 8).Value = 6881.9443   # H132: 5789.636 -> 6881.9443
$ws.Cells.Item(132, 9).Value = 8817   # I132: 6488.6665 -> 8817
$ws.Cells.Item(132, 10).Value = 5914.4165   # J132: 5305.6924 -> 5914.4165
$ws.Cells.Item(132, 11).Value = 26451   # K132: 19465.9995 -> 26451
$ws.Cells.Item(132, 12).Value = 17743.2495   # L132: 15917.0772 -> 17743.2495
$ws.Cells.Item(132, 13).Value = -23921   # M132: -16935.9995 -> -23921
$ws.Cells.Item(132, 14).Value = -22803.2495   # N132: -20977.0772 -> -22803.2495

$ws.Cells.Item(136, 8).Value = 3099.8096   # H136: 2626 -> 3099.8096
$ws.Cells.Item(136, 9).Value = 2392.6428   # I136: 1965.15 -> 2392.6428
$ws.Cells.Item(136, 11).Value = 7177.928400000001   # K136: 5895.450000000001 -> 7177.928400000001
$ws.Cells.Item(136, 13).Value = -4627.928400000001   # M136: -3345.450000000001 -> -4627.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 45304.4   # H16: 44988.4 -> 45304.4
$ws.Cells.Item(16, 10).Value = 45304.4   # J16: 44988.4 -> 45304.4
$ws.Cells.Item(16, 12).Value = 45304.4   # L16: 44988.4 -> 45304.4
$ws.Cells.Item(16, 14).Value = -45888.4   # N16: -45572.4 -> -45888.4

$ws.Cells.Item(113, 8).Value = 743.1111   # H113: 794.2222 -> 743.1111
$ws.Cells.Item(113, 9).Value = 723.5   # I113: 768.5 -> 723.5
$ws.Cells.Item(113, 10).Value = 900   # J113: 1000 -> 900
$ws.Cells.Item(113, 11).Value = 2170.5   # K113: 2305.5 -> 2170.5
$ws.Cells.Item(113, 12).Value = 2700   # L113: 3000 -> 2700
$ws.Cells.Item(113, 13).Value = -0.5   # M113: -135.5 -> -0.5
$ws.Cells.Item(113, 14).Value = -7040   # N113: -7340 -> -7040

$ws.Cells.Item(126, 8).Value = 1400951   # H126: 1470981 -> 1400951
$ws.Cells.Item(126, 9).Value = 1400951   # I126: 1470981 -> 1400951
$ws.Cells.Item(126, 11).Value = 4202853   # K126: 4412943 -> 4202853
$ws.Cells.Item(126, 13).Value = -4200383   # M126: -4410473 -> -4200383

$ws.Cells.Item(135, 8).Value = 47336.117   # H135: 47449.938 -> 47336.117
$ws.Cells.Item(135, 10).Value = 47336.117   # J135: 47449.938 -> 47336.117
$ws.Cells.Item(135, 12).Value = 47336.117   # L135: 47449.938 -> 47336.117
$ws.Cells.Item(135, 14).Value = -57476.117   # N135: -57589.938 -> -57476.117

$ws.Cells.Item(136, 8).Value = 18453.303   # H136: 21116.578 -> 18453.303
$ws.Cells.Item(136, 9).Value = 58772.05   # I136: 65629.12 -> 58772.05
$ws.Cells.Item(136, 10).Value = 2154.2341   # J136: 2198.75 -> 2154.2341
$ws.Cells.Item(136, 11).Value = 176316.15   # K136: 196887.36 -> 176316.15
$ws.Cells.Item(136, 12).Value = 6462.702300000001   # L136: 6596.25 -> 6462.702300000001
$ws.Cells.Item(136, 13).Value = -173766.15   # M136: -194337.36 -> -173766.15
$ws.Cells.Item(136, 14).Value = -11562.7023   # N136: -11696.25 -> -11562.7023

$ws.Cells.Item(137, 8).Value = 49149.734   # H137: 56775.25 -> 49149.734
$ws.Cells.Item(137, 10).Value = 49149.734   # J137: 56775.25 -> 49149.734
$ws.Cells.Item(137, 12).Value = 49149.734   # L137: 56775.25 -> 49149.734
$ws.Cells.Item(137, 14).Value = -59349.734   # N137: -66975.25 -> -59349.734

$ws.Cells.Item(141, 8).Value = 36428.832   # H141: 46857.5 -> 36428.832
$ws.Cells.Item(141, 10).Value = 36428.832   # J141: 46857.5 -> 36428.832
$ws.Cells.Item(141, 12).Value = 36428.832   # L141: 46857.5 -> 36428.832
$ws.Cells.Item(141, 14).Value = -46788.832   # N141: -57217.5 -> -46788.832
